$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 10949
    3  = 13470
    4  = 13939
    5  = 10560
    6  = 12942
    7  = 17625
    8  = 5050
    9  = 14611
    10 = 15741
    11 = 10299
    12 = 17130
    13 = 18324
}

foreach ($row in $values.Keys) {
    $ws.Range("E$row").Value = $values[$row]
}
